# Append the new daily profit row (row 66) to Sheet1, mirroring the format
# of the existing data rows (A = date string, B:L = numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 66

# Keep the date in column A as plain text (matching the other rows) instead
# of letting Excel auto-convert the "mm/dd/yyyy"-looking string into a date
# serial number / date-formatted cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value  = "01/29/2026"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value  = 11269
$ws.Cells.Item($row, 3).Value  = 0.2458624122915205
$ws.Cells.Item($row, 4).Value  = 0.7541375877084795
$ws.Cells.Item($row, 5).Value  = -223.49
$ws.Cells.Item($row, 6).Value  = -30.41
$ws.Cells.Item($row, 7).Value  = -22293.38
$ws.Cells.Item($row, 8).Value  = -72.40000000000001
$ws.Cells.Item($row, 9).Value  = -482.23
$ws.Cells.Item($row, 10).Value = -14.82
$ws.Cells.Item($row, 11).Value = -22775.61
$ws.Cells.Item($row, 12).Value = -66.90000000000001
